# Layout: compacta noticias (tipografia y espaciado)
# Inserts a new slide (duplicate of slide 1, heavily trimmed) as the new
# slide 2, pushing the former slide 2 down to slide 3.

$p = $ppt.ActivePresentation

# --- 1. Duplicate slide 1; PowerPoint inserts the copy right after slide 1,
#        i.e. at index 2, matching the target sldId order (256, 258, 257).
$dupRange = $p.Slides.Item(1).Duplicate()
$ns = $dupRange.Item(1)

# --- 2. Remove the shapes that were deleted from the duplicated slide.
#        (keeps shape id 15 "Imagen 14" and id 35 "Flecha: a la derecha 34")
$idsToDelete = @(19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 30, 32, 34, 36)
foreach ($id in $idsToDelete) {
    for ($i = $ns.Shapes.Count; $i -ge 1; $i--) {
        $sh = $ns.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            $sh.Delete()
        }
    }
}

# --- 3. Reposition / resize the picture that was kept (id 15).
$pic15 = $null
for ($i = 1; $i -le $ns.Shapes.Count; $i++) {
    if ($ns.Shapes.Item($i).Id -eq 15) { $pic15 = $ns.Shapes.Item($i) }
}
$pic15.Left = 5.259212598425197
$pic15.Top = 32.325669291338585
$pic15.Width = 402.11267716535434
$pic15.Height = 335.83929133858265

# The arrow (id 35) keeps the exact same geometry as on slide 1, so no
# change is required there.

# --- 4. Add the three new pictures by duplicating the kept picture and
#        re-cropping / repositioning each copy.

# Imagen 1 (new shape id 2) - plain, uncropped copy of the picture, shown
# to the right of the arrow.
$dup1 = $pic15.Duplicate()
$pic2 = $dup1.Item(1)
$pic2.Name = "Imagen 1"
$pic2.Left = 506.96456692913387
$pic2.Top = 32.325669291338585
$pic2.Width = 402.1127559055118
$pic2.Height = 335.83937007874016

# Imagen 2 (new shape id 3) - cropped detail.
$dup2 = $pic15.Duplicate()
$pic3 = $dup2.Item(1)
$pic3.Name = "Imagen 2"
$pic3.PictureFormat.CropLeft = 10.44384
$pic3.PictureFormat.CropTop = 462.103635
$pic3.PictureFormat.CropRight = 358.76964
$pic3.PictureFormat.CropBottom = 3.9338249999999997
$pic3.Left = 721.4845669291338
$pic3.Top = 202.99732283464567
$pic3.Width = 182.06228346456692
$pic3.Height = 75.00047244094488

# Imagen 3 (new shape id 4) - cropped detail.
$dup3 = $pic15.Duplicate()
$pic4 = $dup3.Item(1)
$pic4.Name = "Imagen 3"
$pic4.PictureFormat.CropLeft = 396.1612575
$pic4.PictureFormat.CropTop = 453.30054
$pic4.PictureFormat.CropRight = 0
$pic4.PictureFormat.CropBottom = 0
$pic4.Left = 519.4086614173228
$pic4.Top = 278.0664566929134
$pic4.Width = 187.34984251968504
$pic4.Height = 90.09850393700788
